# This script applies a cyclic re-shuffle of the per-row data (columns D and L:T)
# across rows 2-11 of the active sheet, leaving row 3 unchanged.
#
# The cycles are (destination row gets the data that used to be in the "next"
# row of the cycle):
#   2 <- 11, 11 <- 4, 4 <- 2
#   5 <- 9,  9 <- 10, 10 <- 5
#   6 <- 8,  8 <- 7,  7 <- 6

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Get-RowData($row) {
    $data = @(
        $ws.Cells.Item($row, 4).Value(),   # D
        $ws.Cells.Item($row, 12).Value(),  # L
        $ws.Cells.Item($row, 13).Value(),  # M
        $ws.Cells.Item($row, 14).Value(),  # N
        $ws.Cells.Item($row, 15).Value(),  # O
        $ws.Cells.Item($row, 16).Value(),  # P
        $ws.Cells.Item($row, 17).Value(),  # Q
        $ws.Cells.Item($row, 18).Value(),  # R
        $ws.Cells.Item($row, 19).Value(),  # S
        $ws.Cells.Item($row, 20).Value()   # T
    )
    return $data
}

function Set-RowData($row, $data) {
    $ws.Cells.Item($row, 4).Value = $data[0]
    $ws.Cells.Item($row, 12).Value = $data[1]
    $ws.Cells.Item($row, 13).Value = $data[2]
    $ws.Cells.Item($row, 14).Value = $data[3]
    $ws.Cells.Item($row, 15).Value = $data[4]
    $ws.Cells.Item($row, 16).Value = $data[5]
    $ws.Cells.Item($row, 17).Value = $data[6]
    $ws.Cells.Item($row, 18).Value = $data[7]
    $ws.Cells.Item($row, 19).Value = $data[8]
    $ws.Cells.Item($row, 20).Value = $data[9]
}

# Snapshot all source rows before writing anything, so cycles do not clobber
# data that is still needed.
$data2  = Get-RowData 2
$data4  = Get-RowData 4
$data5  = Get-RowData 5
$data6  = Get-RowData 6
$data7  = Get-RowData 7
$data8  = Get-RowData 8
$data9  = Get-RowData 9
$data10 = Get-RowData 10
$data11 = Get-RowData 11

# Cycle 1: 2 -> 4 -> 11 -> 2 (row 2 gets old row 11, row 4 gets old row 2, row 11 gets old row 4)
Set-RowData 2  $data11
Set-RowData 4  $data2
Set-RowData 11 $data4

# Cycle 2: 5 -> 10 -> 9 -> 5 (row 5 gets old row 9, row 9 gets old row 10, row 10 gets old row 5)
Set-RowData 5  $data9
Set-RowData 9  $data10
Set-RowData 10 $data5

# Cycle 3: 6 -> 7 -> 8 -> 6 (row 6 gets old row 8, row 7 gets old row 6, row 8 gets old row 7)
Set-RowData 6 $data8
Set-RowData 7 $data6
Set-RowData 8 $data7
